$wb = $excel.ActiveWorkbook

# Column F ("想去人数") updates that apply identically to the "展览" and
# "全部类型" worksheets (rows keyed by row number).
$updates = @{
    2  = 1909
    7  = 1586
    9  = 622
    10 = 365
    19 = 3679
    21 = 4
    23 = 332
    25 = 338
    26 = 345
    28 = 1478
    29 = 143
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
